$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 164, shifting existing rows 164-174 down to 165-175.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record.
$ws.Cells.Item(164, 1).Value = 10
$ws.Cells.Item(164, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(164, 3).Value = "La Araucanía"
$ws.Cells.Item(164, 4).Value = 44753
$ws.Cells.Item(164, 4).NumberFormat = $ws.Cells.Item(165, 4).NumberFormat
$ws.Cells.Item(164, 5).Value = 9
$ws.Cells.Item(164, 6).Value = 100112013
$ws.Cells.Item(164, 7).Value = "Alcachofa"
$ws.Cells.Item(164, 8).Value = "Madrigal"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 240
$ws.Cells.Item(164, 11).Value = 17000
$ws.Cells.Item(164, 12).Value = 18000
$ws.Cells.Item(164, 13).Value = 17500
$ws.Cells.Item(164, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(164, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(164, 16).Value = 438
$ws.Cells.Item(164, 17).Value = 40
$ws.Cells.Item(164, 18).Value = "Hortaliza"
